$d = $word.ActiveDocument

# The first paragraph in the document body holds the hidden ID marker text.
$p1 = $d.Paragraphs(1)

# --- Paragraph formatting -------------------------------------------------
# Add a paragraph border (space-only settings, no visible line) and widen
# the left indent from 120 twips (6 pt) to 225 twips (11.25 pt).
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

$p1.Range.ParagraphFormat.LeftIndent = 11.25

# --- Text content ----------------------------------------------------------
# Replace the marker id and drop the trailing standalone space run, leaving
# a single run with the updated marker text.
$r = $p1.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "**ID__AFFARS_5343_204_70_5__ID**"
